$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 508 (pushes old rows 508-537 down to 510-539)
$ws.Rows("508:509").Insert()

# New row 508 data
$ws.Range("A508").Value = 8
$ws.Range("B508").Value = "Terminal La Palmera de La Serena"
$ws.Range("C508").Value = "Coquimbo"
$ws.Range("D508").Value = 44706
$ws.Range("E508").Value = 4
$ws.Range("F508").Value = 100112043
$ws.Range("G508").Value = "Pepino ensalada"
$ws.Range("H508").Value = "Sin especificar"
$ws.Range("I508").Value = "Primera"
$ws.Range("J508").Value = 700
$ws.Range("K508").Value = 18000
$ws.Range("L508").Value = 19000
$ws.Range("M508").Value = 18500
$ws.Range("N508").Value = "`$/caja 60 unidades"
$ws.Range("O508").Value = "Región de Arica y Parinacota"
$ws.Range("P508").Value = 308
$ws.Range("Q508").Value = 60
$ws.Range("R508").Value = "Hortaliza"

# New row 509 data
$ws.Range("A509").Value = 8
$ws.Range("B509").Value = "Terminal La Palmera de La Serena"
$ws.Range("C509").Value = "Coquimbo"
$ws.Range("D509").Value = 44706
$ws.Range("E509").Value = 4
$ws.Range("F509").Value = 100112043
$ws.Range("G509").Value = "Pepino ensalada"
$ws.Range("H509").Value = "Sin especificar"
$ws.Range("I509").Value = "Segunda"
$ws.Range("J509").Value = 300
$ws.Range("K509").Value = 14000
$ws.Range("L509").Value = 15000
$ws.Range("M509").Value = 14500
$ws.Range("N509").Value = "`$/caja 100 unidades"
$ws.Range("O509").Value = "Región de Arica y Parinacota"
$ws.Range("P509").Value = 145
$ws.Range("Q509").Value = 100
$ws.Range("R509").Value = "Hortaliza"
